$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Enter week 2 (Code Review 2) contribution values for each team member
$ws.Range("D2").Value = 25
$ws.Range("D3").Value = 25
$ws.Range("D4").Value = 25
$ws.Range("D5").Value = 25

# Move selection to C5, matching the author's last click before saving
$ws.Range("C5").Select()
